# Refresh the cryptos price list with the latest scraped values.
# Some "Price" values look like plain numbers (e.g. 487.33); a leading
# apostrophe forces Excel to keep them as text, matching the original
# inline-string cells (and avoiding float round-off such as 487.33 ->
# 487.32999999999998 or the loss of trailing zeros like 1.00 -> 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.412.74'
$ws.Range('E2').Value = '  +1.67%  '
$ws.Range('D3').Value = '3.934.90'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'487.33"
$ws.Range('E5').Value = '  +3.32%  '
$ws.Range('D6').Value = "'148.07"
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('E7').Value = '  +0.82%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'0.735"
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').Value = "'0.169"
$ws.Range('E10').Value = '  +2.77%  '
$ws.Range('E11').Value = '  +4.14%  '
$ws.Range('D12').Value = "'43.04"
$ws.Range('E12').Value = '  -0.90%  '
$ws.Range('D13').Value = "'10.70"
$ws.Range('E13').Value = '  +3.43%  '
$ws.Range('D14').Value = '4.556.58'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.931.43'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').Value = "'14.56"
$ws.Range('E16').Value = '  -3.78%  '
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = "'20.01"
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('E19').Value = '  -1.77%  '
$ws.Range('D20').Value = '68.498.33'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').Value = "'442.84"
$ws.Range('E21').Value = '  +2.48%  '
$ws.Range('E22').Value = '  +4.30%  '
$ws.Range('D23').Value = "'15.15"
$ws.Range('E23').Value = '  +4.06%  '
$ws.Range('D24').Value = "'88.53"
$ws.Range('E24').Value = '  +1.46%  '
$ws.Range('D25').Value = "'11.37"
$ws.Range('E25').Value = '  +17.94%  '
$ws.Range('E26').Value = '  +11.98%  '
$ws.Range('D27').Value = "'3.64"
$ws.Range('E27').Value = '  +1.61%  '
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('D29').Value = "'5.85"
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('D30').Value = "'719.36"
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('D31').Value = "'13.78"
$ws.Range('E31').Value = '  +1.36%  '
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('D33').Value = "'2.91"
$ws.Range('E33').Value = '  +4.22%  '
$ws.Range('D34').Value = "'6.24"
$ws.Range('E34').Value = '  +15.84%  '
$ws.Range('D35').Value = "'42.40"
$ws.Range('E35').Value = '  -0.99%  '
$ws.Range('D36').Value = '0.0₃0892'
$ws.Range('E36').Value = '  +13.69%  '
$ws.Range('D37').Value = "'61.30"
$ws.Range('E37').Value = '  +5.86%  '
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('D39').Value = "'0.399"
$ws.Range('E39').Value = '  +18.86%  '
$ws.Range('D40').Value = "'3.03"
$ws.Range('E40').Value = '  +17.56%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').Value = "'3.24"
$ws.Range('E42').Value = '  +6.53%  '
$ws.Range('D43').Value = "'0.0482"
$ws.Range('E43').Value = '  +1.21%  '
$ws.Range('E44').Value = '  +4.77%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = "'0.143"
$ws.Range('E45').Value = '  +1.15%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = "'1.00"
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0360'
$ws.Range('E47').Value = '  +38.62%  '
$ws.Range('B48').Value = 'LidoDAOToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D48').Value = "'3.42"
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = "'3.25"
$ws.Range('E49').Value = '  +2.32%  '
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('D51').Value = "'145.74"
$ws.Range('E51').Value = '  -0.60%  '
